$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell E1: "DW/cm" (re-uses the existing bold/centred header style)
$ws.Range("E1").Value = "DW/cm"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# New formula columns E (C/B) and F (D/B) for rows 2-23
# Row 2 is entered individually (plain, non-shared formula).
$ws.Range("E2").Formula = "=C2/`$B2"
$ws.Range("F2").Formula = "=D2/`$B2"

# Rows 3-23 are entered as a filled-down range, which the engine stores
# as a shared formula group (matches Excel's own fill-down behaviour).
$ws.Range("E3:E23").Formula = "=C3/`$B3"
$ws.Range("F3:F23").Formula = "=D3/`$B3"

# Update selection to K6
$ws.Range("K6").Select()
